# 25-Jun-2021, end of day update.
# Adds the day's petty-cash transactions (rows 28-33), tops up a couple of
# existing entries (D25, D26, C27) and moves the frozen-pane scroll / active
# selection forward to where the day's work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- existing rows that received extra amounts on this day -----------------

# D25: additional 260,000 on top of the existing 60,000 uang makan
$ws.Range("D25").Formula = "=60000+260000"

# D26: two more purchases added to the running TRANSFER BCA total
$ws.Range("D26").Formula = "=5500000+500000+1000000+1350000+165000+1412500"

# C27: additional 16,258,500 received (TOTAL)
$ws.Range("C27").Formula = "=8000000+16258500"

# --- new transactions recorded for the day ----------------------------------

# Row 28: SALES - cash/retail
$ws.Range("B28").Value = "SALES - cash/retail"
$ws.Range("C28").Formula = "=22885725+1862775-16258500"

# Row 29: SELISIH - lebih
$ws.Range("B29").Value = "SELISIH - lebih"
$ws.Range("C29").Value = 1015000

# Row 30: SETOR KE BANK
$ws.Range("B30").Value = "SETOR KE BANK"
$ws.Range("D30").Value = 23000000

# Row 31: new day (26-Jun-2021 serial 44372) - Wages Expense
$ws.Range("A31").Value = 44372
$ws.Range("B31").Value = "Wages Expense"
$ws.Range("D31").Formula = "=60000"

# Row 32: TRANSFER BCA
$ws.Range("B32").Value = "TRANSFER BCA"
$ws.Range("D32").Formula = "=1740000"

# Row 33: BELi lampu (new expense item / new shared string)
$ws.Range("B33").Value = "BELi lampu"
$ws.Range("D33").Value = 15000

# --- move the view to reflect where the day's entry work ended -------------

$ws.Activate()
$ws.Range("B51").Select()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
